{"js": "// Remove the trailing, now-superfluous single-space run at the end of the\n// list paragraph that ends with \"...folha da \u00e1rvore. \" (note the trailing\n// space after the final period). The paragraph keeps all of its other\n// text/runs untouched; only that last lone-space run disappears.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph by its distinctive, stable text fragment rather than\n// a fixed index, so the script is resilient to unrelated paragraphs being\n// added/removed elsewhere in the document.\nconst marker = \"folha da \u00e1rvore.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error('Could not find the target paragraph containing \"' + marker + '\"');\n}\n\n// Search only within that paragraph's range for single-space matches and\n// take the very last one -- that is the trailing space run to delete.\nconst paragraphRange = target.getRange();\nconst spaceMatches = paragraphRange.search(\" \", { matchCase: true });\nspaceMatches.load(\"items\");\nawait context.sync();\n\nif (spaceMatches.items.length === 0) {\n  throw new Error(\"No trailing space found to remove\");\n}\n\nconst trailingSpace = spaceMatches.items[spaceMatches.items.length - 1];\ntrailingSpace.delete();\nawait context.sync();\n", "ps1": "# Remove the trailing, now-superfluous single-space run at the end of the\n# list paragraph that ends with \"...folha da arvore. \" (note the trailing\n# space after the final period). The rest of the paragraph's text/runs are\n# left untouched; only that last lone-space run disappears.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph by its distinctive, stable text rather than a fixed\n# index, so the script is resilient to unrelated paragraphs being added or\n# removed elsewhere in the document. \"*\" wildcards stand in for the\n# accented characters so the match isn't sensitive to console re-encoding.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*folha da*rvore.*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the target paragraph ending in 'folha da arvore.'\"\n}\n\n# Build a precise range covering only the trailing space: collapse to the\n# end of the paragraph range (which includes the paragraph mark), then walk\n# back over the mark and the space, and drop the mark again from the end.\n$trailing = $target.Range.Duplicate\n$trailing.Collapse(0) | Out-Null       # wdCollapseEnd\n$trailing.MoveStart(1, -2) | Out-Null  # wdCharacter: back over mark + space\n$trailing.MoveEnd(1, -1) | Out-Null    # exclude the paragraph mark again\n\nif ($trailing.Text -ne \" \") {\n    throw \"Unexpected trailing content, expected a single space but found: [$($trailing.Text)]\"\n}\n\n$trailing.Delete()\n"}
